# Swap the full data (columns B..AC) between paired rows in the
# "Sweden Superettan" sheet. The leading index column (A) is left
# untouched; only the row contents are exchanged between each pair.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sweden Superettan")

$firstCol = 2   # B
$lastCol  = 29  # AC

# Two parallel lists (instead of an array-of-pairs) so a single-pair
# workbook edit can never accidentally "unwrap" into a flat list.
$rowsA = @(309, 330, 336, 360, 368, 370, 380, 398, 416, 424, 432, 666, 677, 706, 713)
$rowsB = @(310, 331, 337, 361, 369, 371, 381, 399, 417, 425, 433, 667, 680, 707, 720)

for ($i = 0; $i -lt $rowsA.Count; $i++) {
    $r1 = $rowsA[$i]
    $r2 = $rowsB[$i]

    $range1 = $ws.Range($ws.Cells.Item($r1, $firstCol), $ws.Cells.Item($r1, $lastCol))
    $range2 = $ws.Range($ws.Cells.Item($r2, $firstCol), $ws.Cells.Item($r2, $lastCol))

    $vals1 = $range1.Value2
    $vals2 = $range2.Value2

    $range1.Value2 = $vals2
    $range2.Value2 = $vals1
}
